$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 'D2' '51.848.11'
Set-TextValue 'E2' '  -0.55%  '
Set-TextValue 'D3' '2.811.81'
Set-TextValue 'E3' '  +0.88%  '
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '352.19'
Set-TextValue 'E5' '  +2.28%  '
Set-TextValue 'D6' '111.10'
Set-TextValue 'E6' '  -4.08%  '
Set-TextValue 'E7' '  +2.67%  '
Set-TextValue 'D8' '0.999'
Set-TextValue 'E8' '  -0.03%  '
Set-TextValue 'D9' '0.596'
Set-TextValue 'E9' '  +2.61%  '
Set-TextValue 'D10' '40.48'
Set-TextValue 'E10' '  -5.61%  '
Set-TextValue 'D11' '0.0853'
Set-TextValue 'E11' '  -0.18%  '
Set-TextValue 'E12' '  +0.04%  '
Set-TextValue 'D13' '19.75'
Set-TextValue 'E13' '  -1.84%  '
Set-TextValue 'E14' '  -0.07%  '
Set-TextValue 'D15' '3.248.22'
Set-TextValue 'E15' '  +0.63%  '
Set-TextValue 'D16' '2.810.08'
Set-TextValue 'E16' '  +0.06%  '
Set-TextValue 'D17' '0.918'
Set-TextValue 'E17' '  +3.62%  '
Set-TextValue 'D18' '51.659.28'
Set-TextValue 'E18' '  -0.67%  '
Set-TextValue 'D19' '7.53'
Set-TextValue 'E19' '  +6.81%  '
Set-TextValue 'E20' '  -3.84%  '
Set-TextValue 'E21' '  -0.75%  '
Set-TextValue 'E22' '  +0.83%  '
Set-TextValue 'D23' '70.08'
Set-TextValue 'E23' '  +0.00%  '
Set-TextValue 'D24' '267.76'
Set-TextValue 'E24' '  -0.89%  '
Set-TextValue 'D25' '2.81'
Set-TextValue 'E25' '  +1.49%  '
Set-TextValue 'D26' '26.80'
Set-TextValue 'E26' '  +0.73%  '
Set-TextValue 'E27' '  +0.07%  '
Set-TextValue 'D28' '10.24'
Set-TextValue 'E28' '  -0.32%  '
Set-TextValue 'E29' '  +0.67%  '
Set-TextValue 'E30' '  +20.01%  '
Set-TextValue 'E31' '  -0.42%  '
Set-TextValue 'D32' '52.58'
Set-TextValue 'E32' '  +4.61%  '
Set-TextValue 'D33' '34.18'
Set-TextValue 'E33' '  -1.42%  '
Set-TextValue 'D34' '5.89'
Set-TextValue 'E34' '  +3.08%  '
Set-TextValue 'D35' '5.43'
Set-TextValue 'E35' '  +9.68%  '
Set-TextValue 'E36' '  +2.33%  '
Set-TextValue 'E37' '  -0.11%  '
Set-TextValue 'E38' '  +0.28%  '
Set-TextValue 'D39' '2.02'
Set-TextValue 'E39' '  -4.08%  '
Set-TextValue 'D40' '18.20'
Set-TextValue 'E40' '  -4.40%  '
Set-TextValue 'E41' '  +0.14%  '
Set-TextValue 'D42' '126.86'
Set-TextValue 'E42' '  -0.53%  '
Set-TextValue 'D43' '23.25'
Set-TextValue 'E43' '  -0.82%  '
Set-TextValue 'E45' '  -8.38%  '
Set-TextValue 'D46' '2.087.02'
Set-TextValue 'E46' '  +0.78%  '
Set-TextValue 'E47' '  -0.75%  '
Set-TextValue 'D49' '5.93'
Set-TextValue 'E49' '  +6.84%  '
Set-TextValue 'D50' '0.973'
Set-TextValue 'E50' '  +7.72%  '
Set-TextValue 'D51' '9.02'
Set-TextValue 'E51' '  +1.40%  '
